$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# N6: Last Action Type  "Other" -> "Relationship Building"
$ws.Range("N6").Value = "Relationship Building"

# S6: Next Action Type  "Thank You" -> "Cultivation"
$ws.Range("S6").Value = "Cultivation"

# T6: Next Action Date  45889 -> 45879 (2025-08-20 -> 2025-08-10)
$ws.Range("T6").Value = 45879

# U6: Next Action Assigned To  (empty) -> "Jeff Mcmullen"
$ws.Range("U6").Value = "Jeff Mcmullen"
